$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing style of the D:E data range (rows 2-51), then force
# a text number format so that numeric-looking strings (e.g. "0.9979")
# are written as text instead of being auto-converted to numbers by Excel,
# matching the inlineStr cell type used in the source workbook.
$dataRange = $ws.Range("D2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.319.24'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.864.49'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '0.7135'
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("D6").Value = '238.01'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.9990'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.07941'
$ws.Range("E8").Value = '  -4.09%  '
$ws.Range("D9").Value = '0.3073'
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '25.04'
$ws.Range("E10").Value = '  +7.09%  '
$ws.Range("D11").Value = '0.08166'
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '1.866.64'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '5.234'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '0.7218'
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = '89.29'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '29.342.00'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '5.817'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Value = '241.40'
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("D19").Value = '0.000007817'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = '13.25'
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").Value = '0.9993'
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '2.114.57'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = '0.9978'
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").Value = '7.599'
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("D25").Value = '162.56'
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1457'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '8.933'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '1.917'
$ws.Range("E29").Value = '  -4.19%  '
$ws.Range("D30").Value = '1.373'
$ws.Range("E30").Value = '  -4.51%  '
$ws.Range("D31").Value = '1.471'
$ws.Range("D32").Value = '4.340'
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("D33").Value = '4.049'
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").Value = '0.05194'
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("D35").Value = '1.184'
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("D36").Value = '0.7169'
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("D37").Value = '0.9986'
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").Value = '2.672'
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").Value = '0.01857'
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").Value = '2.698'
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").Value = '1.174.21'
$ws.Range("E41").Value = '  +2.47%  '
$ws.Range("D42").Value = '0.9160'
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").Value = '6.006'
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("D44").Value = '71.69'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").Value = '0.4284'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = '0.9991'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = '102.13'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").Value = '0.5338'
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("D49").Value = '1.760'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").Value = '9.202'
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '7.001'
$ws.Range("E51").Value = '  +0.34%  '

# Restore the original style (removes the temporary text format override)
$dataRange.Style = $origStyle
